$d = $word.ActiveDocument

# --- Update stale cached PAGE field results left over from the last
#     print-preview/field-update pass -------------------------------------
#
# Section 10's primary footer shows the cached page number "27" -> "13"
$footerRange = $d.Sections(10).Footers(1).Range
$footerRange.Find.Execute("27", $false, $true, $false, $false, $false, `
    $true, 1, $false, "13", 2)

# Section 14's primary header shows the cached page number "109" -> "127"
$headerRange = $d.Sections(14).Headers(1).Range
$headerRange.Find.Execute("109", $false, $true, $false, $false, $false, `
    $true, 1, $false, "127", 2)
